$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the reviewer comment anchored on AA7 (it belongs to the data row
# being removed below), then delete that entire row - the "partial match"
# VIN test row (VOLKSWAGEN GOLF / 7MSRP15H&V) - shifting row 8 up to
# become the new (blank) row 7.
$ws.Range("AA7").Comment.Delete()
$ws.Rows.Item(7).Delete()

# Match the resulting active selection on the now-last (blank) row.
$ws.Range("A7:XFD7").Select() | Out-Null
